# Ajout condition nom de sheet
# Rename the original sheet, add a second sheet after it, and update selections.

$wb = $excel.ActiveWorkbook

# Rename "Feuil1" -> "INFOS-ELEVES"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "INFOS-ELEVES"

# Add a new empty sheet "Feuil2" right after INFOS-ELEVES
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Feuil2"

# Set the selection on the new sheet
$ws2.Range("A2:F2").Select()

# Go back to the first sheet and move the selection there too
$ws1.Select()
$ws1.Range("D14").Select()
